$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(1)

# rc5 (id=5)
$sh = $grp.GroupItems.Item(3)
$sh.Left = 179.45363
$sh.Width = 256.8169
$sh.Height = 50.0388

# pl6 (id=6)
$sh = $grp.GroupItems.Item(4)
$sh.Left = 212.9352
$sh.Height = 50.0388

# pl7 (id=7)
$sh = $grp.GroupItems.Item(5)
$sh.Left = 279.8983
$sh.Height = 50.0388

# pl8 (id=8)
$sh = $grp.GroupItems.Item(6)
$sh.Left = 346.8614
$sh.Height = 50.0388

# pl9 (id=9)
$sh = $grp.GroupItems.Item(7)
$sh.Left = 413.8245
$sh.Height = 50.0388

# pl10 (id=10)
$sh = $grp.GroupItems.Item(8)
$sh.Left = 179.45363
$sh.Top = 212.811
$sh.Width = 256.8169

# pl11 (id=11)
$sh = $grp.GroupItems.Item(9)
$sh.Left = 179.45363
$sh.Top = 197.1738
$sh.Width = 256.8169

# pl12 (id=12)
$sh = $grp.GroupItems.Item(10)
$sh.Left = 179.45363
$sh.Top = 181.5367
$sh.Width = 256.8169

# pl13 (id=13)
$sh = $grp.GroupItems.Item(11)
$sh.Left = 179.45363
$sh.Height = 50.0388

# pl14 (id=14)
$sh = $grp.GroupItems.Item(12)
$sh.Left = 246.4167
$sh.Height = 50.0388

# pl15 (id=15)
$sh = $grp.GroupItems.Item(13)
$sh.Left = 313.3799
$sh.Height = 50.0388

# pl16 (id=16)
$sh = $grp.GroupItems.Item(14)
$sh.Left = 380.343
$sh.Height = 50.0388

# rc17 (id=17)
$sh = $grp.GroupItems.Item(15)
$sh.Left = 179.45363
$sh.Top = 174.50
$sh.Width = 251.7813
$sh.Height = 14.0734

# rc18 (id=18)
$sh = $grp.GroupItems.Item(16)
$sh.Left = 179.45363
$sh.Top = 190.1371
$sh.Width = 18.08001
$sh.Height = 14.0734

# rc19 (id=19)
$sh = $grp.GroupItems.Item(17)
$sh.Left = 179.45363
$sh.Top = 205.7743
$sh.Width = 6.8303
$sh.Height = 14.0734

# pg20 (id=20)
$sh = $grp.GroupItems.Item(18)
$sh.Left = 408.9064
$sh.Top = 176.6632

# pg21 (id=21)
$sh = $grp.GroupItems.Item(19)
$sh.Left = 410.2401
$sh.Top = 177.4301

# pg22 (id=22)
$sh = $grp.GroupItems.Item(20)
$sh.Left = 417.2422
$sh.Top = 176.6966

# pg23 (id=23)
$sh = $grp.GroupItems.Item(21)
$sh.Left = 420.10961
$sh.Top = 176.3965

# pg24 (id=24)
$sh = $grp.GroupItems.Item(22)
$sh.Left = 420.7099
$sh.Top = 178.8305

# pg25 (id=25)
$sh = $grp.GroupItems.Item(23)
$sh.Left = 419.7429
$sh.Top = 180.331

# pg26 (id=26)
$sh = $grp.GroupItems.Item(24)
$sh.Left = 204.0909
$sh.Top = 193.24142

# pg27 (id=27)
$sh = $grp.GroupItems.Item(25)
$sh.Left = 206.3582
$sh.Top = 194.8419

# pg28 (id=28)
$sh = $grp.GroupItems.Item(26)
$sh.Left = 206.3582
$sh.Top = 196.209

# pg29 (id=29)
$sh = $grp.GroupItems.Item(27)
$sh.Left = 206.3582
$sh.Top = 197.6427

# pg30 (id=30)
$sh = $grp.GroupItems.Item(28)
$sh.Left = 204.2243
$sh.Top = 199.2432

# pg31 (id=31)
$sh = $grp.GroupItems.Item(29)
$sh.Left = 208.6589
$sh.Top = 199.2432

# pg32 (id=32)
$sh = $grp.GroupItems.Item(30)
$sh.Left = 212.8268
$sh.Top = 193.1415

# pg33 (id=33)
$sh = $grp.GroupItems.Item(31)
$sh.Left = 214.2272
$sh.Top = 195.6422

# pg34 (id=34)
$sh = $grp.GroupItems.Item(32)
$sh.Left = 192.94111
$sh.Top = 208.8452

# pg35 (id=35)
$sh = $grp.GroupItems.Item(33)
$sh.Left = 195.842
$sh.Top = 211.01252

# pg36 (id=36)
$sh = $grp.GroupItems.Item(34)
$sh.Left = 197.109
$sh.Top = 213.1798

# pg37 (id=37)
$sh = $grp.GroupItems.Item(35)
$sh.Left = 192.97441
$sh.Top = 210.9125

# pg38 (id=38)
$sh = $grp.GroupItems.Item(36)
$sh.Left = 193.0411
$sh.Top = 215.2804

# pg39 (id=39)
$sh = $grp.GroupItems.Item(37)
$sh.Left = 201.2769
$sh.Top = 208.8119

# pg40 (id=40)
$sh = $grp.GroupItems.Item(38)
$sh.Left = 203.3108
$sh.Top = 210.479

# pg41 (id=41)
$sh = $grp.GroupItems.Item(39)
$sh.Left = 203.9777
$sh.Top = 213.8467

# rc42 (id=42)
$sh = $grp.GroupItems.Item(40)
$sh.Left = 179.45363
$sh.Width = 256.8169
$sh.Height = 50.0388

# pg43 (id=43)
$sh = $grp.GroupItems.Item(41)
$sh.Left = 172.0217
$sh.Top = 209.6936

# pg44 (id=44)
$sh = $grp.GroupItems.Item(42)
$sh.Left = 172.12481
$sh.Top = 194.0565

# pg45 (id=45)
$sh = $grp.GroupItems.Item(43)
$sh.Left = 172.37402
$sh.Top = 178.5096

# pl46 (id=46)
$sh = $grp.GroupItems.Item(44)
$sh.Left = 176.7139
$sh.Top = 212.811

# pl47 (id=47)
$sh = $grp.GroupItems.Item(45)
$sh.Left = 176.7139
$sh.Top = 197.1738

# pl48 (id=48)
$sh = $grp.GroupItems.Item(46)
$sh.Left = 176.7139
$sh.Top = 181.5367

# pl49 (id=49)
$sh = $grp.GroupItems.Item(47)
$sh.Left = 179.45363
$sh.Top = 222.1933

# pl50 (id=50)
$sh = $grp.GroupItems.Item(48)
$sh.Left = 246.4167
$sh.Top = 222.1933

# pl51 (id=51)
$sh = $grp.GroupItems.Item(49)
$sh.Left = 313.3799
$sh.Top = 222.1933

# pl52 (id=52)
$sh = $grp.GroupItems.Item(50)
$sh.Left = 380.343
$sh.Top = 222.1933

# pg53 (id=53)
$sh = $grp.GroupItems.Item(51)
$sh.Left = 177.3503
$sh.Top = 227.0345

# pg54 (id=54)
$sh = $grp.GroupItems.Item(52)
$sh.Left = 178.1323
$sh.Top = 227.66182

# pg55 (id=55)
$sh = $grp.GroupItems.Item(53)
$sh.Left = 239.4278
$sh.Top = 227.1248

# pg56 (id=56)
$sh = $grp.GroupItems.Item(54)
$sh.Left = 244.3134
$sh.Top = 227.0345

# pg57 (id=57)
$sh = $grp.GroupItems.Item(55)
$sh.Left = 245.0955
$sh.Top = 227.66182

# pg58 (id=58)
$sh = $grp.GroupItems.Item(56)
$sh.Left = 249.2076
$sh.Top = 227.0345

# pg59 (id=59)
$sh = $grp.GroupItems.Item(57)
$sh.Left = 249.9896
$sh.Top = 227.66182

# pg60 (id=60)
$sh = $grp.GroupItems.Item(58)
$sh.Left = 304.26183
$sh.Top = 227.1248

# pg61 (id=61)
$sh = $grp.GroupItems.Item(59)
$sh.Left = 308.8295
$sh.Top = 227.0345

# pg62 (id=62)
$sh = $grp.GroupItems.Item(60)
$sh.Left = 309.61142
$sh.Top = 227.66182

# pg63 (id=63)
$sh = $grp.GroupItems.Item(61)
$sh.Left = 313.7236
$sh.Top = 227.0345

# pg64 (id=64)
$sh = $grp.GroupItems.Item(62)
$sh.Left = 314.5056
$sh.Top = 227.66182

# pg65 (id=65)
$sh = $grp.GroupItems.Item(63)
$sh.Left = 318.61773
$sh.Top = 227.0345

# pg66 (id=66)
$sh = $grp.GroupItems.Item(64)
$sh.Left = 319.3998
$sh.Top = 227.66182

# pg67 (id=67)
$sh = $grp.GroupItems.Item(65)
$sh.Left = 371.225
$sh.Top = 227.1248

# pg68 (id=68)
$sh = $grp.GroupItems.Item(66)
$sh.Left = 375.80111
$sh.Top = 227.1248

# pg69 (id=69)
$sh = $grp.GroupItems.Item(67)
$sh.Left = 380.6867
$sh.Top = 227.0345

# pg70 (id=70)
$sh = $grp.GroupItems.Item(68)
$sh.Left = 381.4688
$sh.Top = 227.66182

# pg71 (id=71)
$sh = $grp.GroupItems.Item(69)
$sh.Left = 385.5808
$sh.Top = 227.0345

# pg72 (id=72)
$sh = $grp.GroupItems.Item(70)
$sh.Left = 386.36284
$sh.Top = 227.66182
